# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("1000 Bs = 12.24 = 49555.03 pesos", "1000 Bs = 11.98 = 48455.66 pesos")
$text = $text.Replace("49555.03 pesos = 12.24 = 980.59 Bs", "48455.66 pesos = 11.92 = 969.23 Bs")
$cell.Value = $text

# --- tasas: update rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 83.44
$ws2.Range("O10").Value = 4043.14
$ws2.Range("N12").Value = 4064.5
$ws2.Range("O12").Value = 81.3
